$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 282.16666
$ws.Range("I9").Value = 336.2
$ws.Range("J9").Value = 12
$ws.Range("K9").Value = 336.2
$ws.Range("L9").Value = 12
$ws.Range("M9").Value = -167.2
$ws.Range("N9").Value = -350

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()

$ws.Range("H31").Value = 2863.5
$ws.Range("I31").Value = 45.25
$ws.Range("K31").Value = 135.75
$ws.Range("M31").Value = 94.25

$ws.Range("H94").Value = 9135.223
$ws.Range("I94").Value = 9135.223
$ws.Range("K94").Value = 9135.223
$ws.Range("M94").Value = -8684.223

$ws.Range("H96").Value = 216.33333
$ws.Range("I96").Value = 216.33333
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 648.99999
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 724.00001
$ws.Range("N96").ClearContents()

$ws.Range("H132").Value = 3251.238
$ws.Range("I132").Value = 804.4706
$ws.Range("K132").Value = 2413.4118
$ws.Range("M132").Value = 116.5882000000001

$ws.Range("H137").Value = 2293.348
$ws.Range("I137").Value = 1411.8462
$ws.Range("K137").Value = 4235.5386
$ws.Range("M137").Value = -1685.5386


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6799.125
$ws.Range("I2").Value = 4367.077
$ws.Range("J2").Value = 17338
$ws.Range("K2").Value = 4367.077
$ws.Range("L2").Value = 17338
$ws.Range("M2").Value = -4254.077
$ws.Range("N2").Value = -17564

$ws.Range("H32").Value = 4554879.5
$ws.Range("I32").Value = 5615.75
$ws.Range("J32").Value = 16686250
$ws.Range("K32").Value = 5615.75
$ws.Range("L32").Value = 16686250
$ws.Range("M32").Value = -5328.75
$ws.Range("N32").Value = -16686824

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H74").Value = 7136.647
$ws.Range("I74").Value = 6776.357
$ws.Range("K74").Value = 6776.357
$ws.Range("M74").Value = -5902.357

$ws.Range("H77").Value = 7136.647
$ws.Range("I77").Value = 6776.357
$ws.Range("K77").Value = 33881.785
$ws.Range("M77").Value = -29513.785

$ws.Range("H116").Value = 6799.125
$ws.Range("I116").Value = 4367.077
$ws.Range("J116").Value = 17338
$ws.Range("K116").Value = 4367.077
$ws.Range("L116").Value = 17338
$ws.Range("M116").Value = -2073.077
$ws.Range("N116").Value = -21926

$ws.Range("H122").Value = 2665.6667
$ws.Range("I122").Value = 2665.6667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7997.000100000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5547.000100000001
$ws.Range("N122").ClearContents()


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6799.125
$ws.Range("I3").Value = 4367.077
$ws.Range("J3").Value = 17338
$ws.Range("K3").Value = 4367.077
$ws.Range("L3").Value = 17338
$ws.Range("M3").Value = -4253.077
$ws.Range("N3").Value = -17566

$ws.Range("H86").Value = 6845.615
$ws.Range("I86").Value = 2501.5
$ws.Range("J86").Value = 7635.4546
$ws.Range("K86").Value = 2501.5
$ws.Range("L86").Value = 7635.4546
$ws.Range("M86").Value = -1378.5
$ws.Range("N86").Value = -9881.454600000001

$ws.Range("H89").Value = 6845.615
$ws.Range("I89").Value = 2501.5
$ws.Range("J89").Value = 7635.4546
$ws.Range("K89").Value = 12507.5
$ws.Range("L89").Value = 38177.273
$ws.Range("M89").Value = -6891.5
$ws.Range("N89").Value = -49409.273

$ws.Range("H99").Value = 3649.1667
$ws.Range("I99").Value = 3098
$ws.Range("K99").Value = 3098
$ws.Range("M99").Value = -1600

$ws.Range("H134").Value = 1800.12
$ws.Range("I134").Value = 1800.12
$ws.Range("K134").Value = 5400.36
$ws.Range("M134").Value = -2865.36


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4451.457
$ws.Range("J132").Value = 6171.5454
$ws.Range("L132").Value = 18514.6362
$ws.Range("N132").Value = -23574.6362

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 6249.5
$ws.Range("J39").Value = 6681.273
$ws.Range("L39").Value = 20043.819
$ws.Range("N39").Value = -20631.819

$ws.Range("H68").Value = 990.7
$ws.Range("J68").Value = 1108.8
$ws.Range("L68").Value = 3326.4
$ws.Range("N68").Value = -4948.4

$ws.Range("H71").Value = 990.7
$ws.Range("J71").Value = 1108.8
$ws.Range("L71").Value = 9979.199999999999
$ws.Range("N71").Value = -18091.2


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 11653265
$ws.Range("I11").Value = 14666800
$ws.Range("J11").Value = 6002887.5
$ws.Range("K11").Value = 14666800
$ws.Range("L11").Value = 6002887.5
$ws.Range("M11").Value = -14666661
$ws.Range("N11").Value = -6003165.5

$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 874.6667
$ws.Range("I16").Value = 874.6667
$ws.Range("K16").Value = 874.6667
$ws.Range("M16").Value = -704.6667

$ws.Range("H46").Value = 6132.864
$ws.Range("J46").Value = 6187.1875
$ws.Range("L46").Value = 6187.1875
$ws.Range("N46").Value = -6563.1875

$ws.Range("H55").Value = 866
$ws.Range("I55").Value = 1064
$ws.Range("J55").Value = 470
$ws.Range("K55").Value = 1064
$ws.Range("L55").Value = 470
$ws.Range("M55").Value = -891
$ws.Range("N55").Value = -816

$ws.Range("H122").Value = 8499.5
$ws.Range("I122").Value = 7000
$ws.Range("J122").Value = 9999
$ws.Range("K122").Value = 21000
$ws.Range("L122").Value = 29997
$ws.Range("M122").Value = -18550
$ws.Range("N122").Value = -34897


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 182423.08
$ws.Range("I4").Value = 182423.08
$ws.Range("K4").Value = 182423.08
$ws.Range("M4").Value = -182310.08

